$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D (rows 2-6) with the new "U" label
$ws.Range("D2:D6").Value = "U"

# Remove the now-unused sire/dam columns (F, G) for data rows
$ws.Range("F2:G6").ClearContents()

# Leave the selection where the author left it after editing
$ws.Range("D7").Select()
